# no-op for now
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
